$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Append the new testmail row (row 15) to the Logs sheet ---
$ws.Range("A15").Value = "Wil je deze klant bellen?"
$ws.Range("B15").Value = "mailmind.test@zohomail.eu"
$ws.Range("C15").Value = "Testmail #5: Wil je deze klant bellen?"
$ws.Range("D15").Value = "Overig"
$ws.Range("E15").Value = "Geachte klantenservice,`nDit is slechts een testmail om te controleren of het systeem werkt. Er is geen echte klant die gebeld hoeft te worden.`nMet vriendelijke groet,`n[Naam]"
$ws.Range("F15").Value = "2025-08-03 14:42:28"
$ws.Range("G15").Value = "Ja"
$ws.Range("H15").Value = "Nee"
$ws.Range("I15").Value = "Ja"
$ws.Range("J15").Value = "Nee"

# --- Extend the existing conditional-formatting rules to cover the new row ---
$colRanges = @("D2:D14", "G2:G14", "H2:H14", "I2:I14", "J2:J14")
foreach ($colRange in $colRanges) {
  $col = $colRange.Substring(0, 1)
  $newRange = $ws.Range($col + "2:" + $col + "15")
  $fcs = $ws.Range($colRange).FormatConditions
  for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($newRange)
  }
}

# --- Update the Dashboard summary count for the "Overig" category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B4").Value = 4
